$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corrections (Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020))
$ws.Range("H3").Value = 0.3852459016393442
$ws.Range("I3").Value = 0.1636363636363636
$ws.Range("K3").Value = 87.3

$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 25
$ws.Range("S3").Value = 37
$ws.Range("T3").Value = 51
$ws.Range("U3").Value = 108
$ws.Range("V3").Value = 280
$ws.Range("W3").Value = 261
$ws.Range("X3").Value = 249
$ws.Range("Y3").Value = 235
$ws.Range("Z3").Value = 178

$ws.Range("AF3").Value = 0.979021
$ws.Range("AG3").Value = 0.912587
$ws.Range("AH3").Value = 0.870629
$ws.Range("AI3").Value = 0.821678
$ws.Range("AJ3").Value = 0.622378
